$d = $word.ActiveDocument

# Locate the existing "Quantity?" bullet (numId=3, ilvl=1) and position a
# range right after it, so the new "Date posted" bullet is inserted as the
# next sibling paragraph with the same list/indent formatting.
$r = $d.Content
$r.Find.Execute("Quantity?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Collapse to the end of "Quantity?", insert a new paragraph mark there
# (inherits the paragraph/run formatting of "Quantity?"), then move past
# the new paragraph mark and type the new bullet's text.
$r.Collapse(0)
$r.InsertParagraphAfter()
$r.Move(1, 1)
$r.InsertAfter("Date posted")
